$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: convert from text "76442711" to a numeric value 76442711
$ws.Range("A2").Value = 76442711

# C2: update points value from 0 to 408
$ws.Range("C2").Value = 408
